# Update all final evaluation results across the three sheets of the
# evaluation_metrics workbook (Summary, Classification Report, Confusion Matrix).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Summary"
# Columns: A=Model, B=Accuracy, C=Precision, D=Recall, E=F1 Score,
#          F=F2 Score, G=F5 Score, H=AUC, I=True Positives, J=False Positives,
#          K=True Negatives, L=False Negatives
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 0.0498220640569395
$wsSummary.Range("C2").Value = 0.0498220640569395
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.09491525423728814
$wsSummary.Range("F2").Value = 0.2077151335311573
$wsSummary.Range("G2").Value = 0.5768621236133122
$wsSummary.Range("H2").Value = 0.5801899411449973
$wsSummary.Range("I2").Value = 28
$wsSummary.Range("J2").Value = 534
$wsSummary.Range("K2").Value = 0
$wsSummary.Range("L2").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Classification Report"
# Columns: B=precision, C=recall, D=f1-score, E=support
# Rows: 2="0", 3="1", 4="accuracy", 5="macro avg", 6="weighted avg"
# ---------------------------------------------------------------------------
$wsReport = $wb.Worksheets.Item("Classification Report")

# row 2 - label "0"
$wsReport.Range("B2").Value = 0
$wsReport.Range("C2").Value = 0
$wsReport.Range("D2").Value = 0

# row 3 - label "1"
$wsReport.Range("B3").Value = 0.0498220640569395
$wsReport.Range("C3").Value = 1
$wsReport.Range("D3").Value = 0.09491525423728814

# row 4 - "accuracy"
$wsReport.Range("B4").Value = 0.0498220640569395
$wsReport.Range("C4").Value = 0.0498220640569395
$wsReport.Range("D4").Value = 0.0498220640569395
$wsReport.Range("E4").Value = 0.0498220640569395

# row 5 - "macro avg"
$wsReport.Range("B5").Value = 0.02491103202846975
$wsReport.Range("D5").Value = 0.04745762711864407

# row 6 - "weighted avg"
$wsReport.Range("B6").Value = 0.002482238066893783
$wsReport.Range("C6").Value = 0.0498220640569395
$wsReport.Range("D6").Value = 0.004728873876590867

# ---------------------------------------------------------------------------
# Sheet 3: "Confusion Matrix"
# Columns: B=Predicted 0, C=Predicted 1
# Rows: 2="Actual 0", 3="Actual 1"
# ---------------------------------------------------------------------------
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")

$wsConfusion.Range("B2").Value = 0
$wsConfusion.Range("C2").Value = 534

$wsConfusion.Range("B3").Value = 0
$wsConfusion.Range("C3").Value = 28
